$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 value: "Empresa Id" -> "Entidad Id"
$ws.Range("B1").Value = "Entidad Id"

# Update selection to D3
$ws.Range("D3").Select()

# Update the window size/position seen in the bookViews (best-effort; this
# reflects the native Excel application window chrome, not document content).
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 20490
$excel.ActiveWindow.Height = 10920
